# Add data for 2022-03-31
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-23"

# Update the column header label (shared string "2022 (through 03-22)")
$ws.Range("I1").Value = "2022 (through 03-23)"

# Update the March 2022 value and the running Total for that column
$ws.Range("I4").Value = 97
$ws.Range("I14").Value = 397
